$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.875.24"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.995.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.51%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.22"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.71%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.71"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.95%  "

$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.477.99"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.997.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.891.71"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.41"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.77"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.177"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.77"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("E30").Value = "  +1.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.12"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.83%  "

$ws.Range("E33").Value = "  +12.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.62"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0431"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.16"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.73%  "

$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.119.83"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.41%  "

$ws.Range("E48").Value = "  -7.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.240"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0333"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.895"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.75%  "
